$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '26.908.97'
Set-TextValue $ws.Range('E2') '  -0.23%  '
Set-TextValue $ws.Range('D3') '1.814.24'
Set-TextValue $ws.Range('E3') '  +0.17%  '
Set-TextValue $ws.Range('D5') '309.28'
Set-TextValue $ws.Range('E5') '  -0.45%  '
Set-TextValue $ws.Range('E6') '  +0.12%  '
Set-TextValue $ws.Range('D7') '0.4660'
Set-TextValue $ws.Range('E7') '  +0.68%  '
Set-TextValue $ws.Range('D8') '0.3654'
Set-TextValue $ws.Range('E8') '  -1.75%  '
Set-TextValue $ws.Range('E9') '  -0.27%  '
Set-TextValue $ws.Range('D10') '0.8690'
Set-TextValue $ws.Range('E10') '  -0.76%  '
Set-TextValue $ws.Range('D11') '20.22'
Set-TextValue $ws.Range('E11') '  -1.28%  '
Set-TextValue $ws.Range('D12') '1.839.59'
Set-TextValue $ws.Range('E12') '  +1.47%  '
Set-TextValue $ws.Range('D13') '5.375'
Set-TextValue $ws.Range('D14') '0.07097'
Set-TextValue $ws.Range('E14') '  +0.47%  '
Set-TextValue $ws.Range('D15') '6.498'
Set-TextValue $ws.Range('E15') '  -0.56%  '
Set-TextValue $ws.Range('D16') '91.12'
Set-TextValue $ws.Range('E16') '  -1.55%  '
Set-TextValue $ws.Range('E17') '  +0.19%  '
Set-TextValue $ws.Range('D18') '0.000008688'
Set-TextValue $ws.Range('E18') '  -0.33%  '
Set-TextValue $ws.Range('E19') '  +0.10%  '
Set-TextValue $ws.Range('E20') '  -0.90%  '
Set-TextValue $ws.Range('D21') '26.927.50'
Set-TextValue $ws.Range('E21') '  -0.14%  '
Set-TextValue $ws.Range('D22') '5.289'
Set-TextValue $ws.Range('E22') '  -0.54%  '
Set-TextValue $ws.Range('E23') '  -0.77%  '
Set-TextValue $ws.Range('D24') '2.046.94'
Set-TextValue $ws.Range('E24') '  +0.35%  '
Set-TextValue $ws.Range('D25') '1.894'
Set-TextValue $ws.Range('E25') '  -0.21%  '
Set-TextValue $ws.Range('D26') '150.81'
Set-TextValue $ws.Range('E26') '  -0.50%  '
Set-TextValue $ws.Range('D27') '18.31'
Set-TextValue $ws.Range('E27') '  -0.16%  '
Set-TextValue $ws.Range('D28') '2.131'
Set-TextValue $ws.Range('E28') '  -1.07%  '
Set-TextValue $ws.Range('D29') '5.245'
Set-TextValue $ws.Range('E29') '  -0.82%  '
Set-TextValue $ws.Range('D30') '115.64'
Set-TextValue $ws.Range('E30') '  -0.25%  '
Set-TextValue $ws.Range('D31') '0.08911'
Set-TextValue $ws.Range('E31') '  -0.28%  '
Set-TextValue $ws.Range('D32') '0.7536'
Set-TextValue $ws.Range('E32') '  -0.29%  '
Set-TextValue $ws.Range('D33') '1.162'
Set-TextValue $ws.Range('E33') '  +0.30%  '
Set-TextValue $ws.Range('E34') '  +0.29%  '
Set-TextValue $ws.Range('D35') '2.907'
Set-TextValue $ws.Range('E35') '  -0.56%  '
Set-TextValue $ws.Range('E36') '  +0.15%  '
Set-TextValue $ws.Range('D37') '1.092'
Set-TextValue $ws.Range('E37') '  -1.32%  '
Set-TextValue $ws.Range('E38') '  +0.55%  '
Set-TextValue $ws.Range('E39') '  -1.83%  '
Set-TextValue $ws.Range('D40') '2.969'
Set-TextValue $ws.Range('E40') '  +1.58%  '
Set-TextValue $ws.Range('D41') '7.194'
Set-TextValue $ws.Range('E41') '  -0.17%  '
Set-TextValue $ws.Range('D42') '0.5280'
Set-TextValue $ws.Range('E42') '  -0.78%  '
Set-TextValue $ws.Range('D43') '2.305'
Set-TextValue $ws.Range('E43') '  -6.21%  '
Set-TextValue $ws.Range('D44') '0.1655'
Set-TextValue $ws.Range('E44') '  -0.61%  '
Set-TextValue $ws.Range('D45') '8.415'
Set-TextValue $ws.Range('E45') '  -1.41%  '
Set-TextValue $ws.Range('D46') '0.4842'
Set-TextValue $ws.Range('E46') '  -3.04%  '
Set-TextValue $ws.Range('D47') '10.44'
Set-TextValue $ws.Range('E47') '  +0.77%  '
Set-TextValue $ws.Range('E48') '  +0.16%  '
Set-TextValue $ws.Range('D49') '103.11'
Set-TextValue $ws.Range('E49') '  -1.15%  '
Set-TextValue $ws.Range('E50') '  -1.01%  '
